$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = "-"

$ws.Range("B3").Value = "[-, -, 'MCT-1A-Metrologia', -]"
$ws.Range("C3").Value = "-"

$ws.Range("B4").Value = "[-, -, 'MCT-1A-Metrologia', -]"
$ws.Range("C4").Value = "-"

$ws.Range("B6").Value = "[-, -, 'MCT-1A-Metrologia', -]"
$ws.Range("C6").Value = "-"

$ws.Range("B7").Value = "[-, -, 'MCT-1A-Metrologia', -]"
$ws.Range("C7").Value = "-"

$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = "-"

$ws.Range("E10").Value = "-"

$ws.Range("E11").Value = "-"

$ws.Range("E12").Value = "-"

$ws.Range("E16").Value = "-"

$ws.Range("B18").Value = "-"

$ws.Range("B19").Value = "-"

$ws.Range("B20").Value = "-"
$ws.Range("E20").Value = "-"
$ws.Range("F20").Value = "-"

$ws.Range("B21").Value = "-"
$ws.Range("E21").Value = "-"
$ws.Range("F21").Value = "-"
